# "Generate Report for handback" — mark the pending rows as handed back,
# and record the resulting target/handback file links + handback timestamp.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$hlColor = 15570276   # RGB(0x64,0x95,0xED) -> matches the workbook's existing hyperlink font color

function Apply-HyperlinkLook($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hlColor
}

# ---------------------------------------------------------------------------
# Overview sheet: just the status text changes (columns B/C mirror zh-cn/de-de)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
foreach ($addr in @("B2", "C2", "B3", "C3")) {
    $cell = $overview.Range($addr)
    if ($cell.Value2 -eq $statusOld) {
        $cell.Value = $statusNew
    }
}

# ---------------------------------------------------------------------------
# Per-language detail sheets: zh-cn / de-de
# ---------------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; HandbackDateTime = "2016-01-26 03:42:33"; HandbackRepoBase = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e74167aee4df4786cab22121003addbd331ae157/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/" },
    @{ Name = "de-de"; HandbackDateTime = "2016-01-26 03:42:49"; HandbackRepoBase = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/631abadaaddda65167dfb8b57eb5aedb8effef4b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/" }
)

$sourceRepoBase = "https://github.com/OpenLocalizationTest/oltest/blob/9109eec8c405cb6109142aeeba91d9840004b94a/e2e/"

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    foreach ($row in @(2, 3)) {
        $aCell = $ws.Range("A$row")
        $cCell = $ws.Range("C$row")
        $bCell = $ws.Range("B$row")

        $sourceFileName = $aCell.Value2
        $handoffFileName = $cCell.Value2

        # Only rows that were "Ready for handoff" get handed back — the
        # ".localization-config" / "Ignored" row (row 4) is left untouched.
        if ($bCell.Value2 -ne $statusOld) {
            continue
        }

        # Status -> handed back, in sync with en-US
        $bCell.Value = $statusNew

        # Latest Target File (E): same file name as the source markdown file
        $eCell = $ws.Range("E$row")
        $eCell.Value = $sourceFileName
        $ws.Hyperlinks.Add($eCell, ($sourceRepoBase + $sourceFileName), [Type]::Missing, [Type]::Missing, $sourceFileName) | Out-Null
        Apply-HyperlinkLook $eCell

        # Latest Handback File (F): same file name as the latest handoff xlf
        $fCell = $ws.Range("F$row")
        $fCell.Value = $handoffFileName
        $ws.Hyperlinks.Add($fCell, ($lang.HandbackRepoBase + $handoffFileName), [Type]::Missing, [Type]::Missing, $handoffFileName) | Out-Null
        Apply-HyperlinkLook $fCell

        # Latest Handback DateTime (G): stamp with the handback time
        $gCell = $ws.Range("G$row")
        $gCell.Value = $lang.HandbackDateTime
    }
}

"Handback report generated."
